$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 0) Grab a formatted copy of the "Meta description" paragraph's
#    run structure (empty run + bold run) before it gets removed;
#    we reuse its exact run layout for the new bold paragraph that
#    gets added near the end of the document.
# ---------------------------------------------------------------
$metaRange = $d.Paragraphs(2).Range
$metaFormattedText = $metaRange.FormattedText

# ---------------------------------------------------------------
# 1) Page title (Heading1)
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Treasures of the Dead for Free: Game Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Treasures of the Dead Free - Exciting Slot Game", 2) | Out-Null

# ---------------------------------------------------------------
# 2) "What we like" bullet list
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Interesting gameplay mechanics with Hyperlines and Bonus Game", $true, $false, $false, $false, $false,
    $true, 1, $false, "Easy Autoplay feature for convenient gameplay", 2) | Out-Null

$d.Content.Find.Execute(
    "Autoplay and Turbo features for convenient gameplay", $true, $false, $false, $false, $false,
    $true, 1, $false, "Turbo button options to speed up gameplay", 2) | Out-Null

$d.Content.Find.Execute(
    "Well-designed graphics and high-quality sound", $true, $false, $false, $false, $false,
    $true, 1, $false, "Well-designed graphics with large symbols", 2) | Out-Null

$d.Content.Find.Execute(
    "Mobile compatibility for playing on-the-go", $true, $false, $false, $false, $false,
    $true, 1, $false, "Rock version of eastern sounds adds excitement", 2) | Out-Null

# ---------------------------------------------------------------
# 3) "What we don't like" bullet list
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Limited number of paylines with fixed values", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited number of fixed paylines", 2) | Out-Null

$d.Content.Find.Execute(
    "Bonus Game can only be accessed by purchasing it, limiting accessibility", $true, $false, $false, $false, $false,
    $true, 1, $false, "Expanding special symbol only available during Bonus Game", 2) | Out-Null

# ---------------------------------------------------------------
# 4) Insert a new bold paragraph right after the "What we don't
#    like" list (before the final italic image-prompt paragraph),
#    reusing the run layout captured in step 0.
# ---------------------------------------------------------------
$lastBullet = $d.Content.Find
$anchor = $d.Content
$anchor.Find.Execute("Expanding special symbol only available during Bonus Game") | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()

# Re-locate the freshly inserted (still empty) paragraph: it is the
# paragraph right after the one containing our anchor text.
$allParas = $d.Paragraphs
$insertIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs($insertIndex)
$newPara.Style = -1   # wdStyleNormal -- clears the inherited ListBullet style
$newPara.Range.FormattedText = $metaFormattedText

$newParaRange = $d.Paragraphs($insertIndex).Range
$newParaRange.Find.Execute(
    "Meta description", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Treasures of the Dead Free - Exciting Slot Game", 2) | Out-Null

$newParaRange2 = $d.Paragraphs($insertIndex).Range
$newParaRange2.Find.Execute(
    ": Explore the features of Treasures of the Dead online slot game in our review. Play for free and enjoy Hyperlines, Bonus Game, and mobile compatibility.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------
# 5) Replace the final italic "feature image prompt" paragraph
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a Feature Image Prompt: Create a feature image for Treasures of the Dead that incorporates a cartoon-style Maya warrior with glasses. The warrior should be depicted as happy and triumphant, holding a treasure chest or a magical book in one hand, while the other hand is raised in victory. The background of the image should showcase the inside of a pyramid or a tomb, with hidden treasures, jewels, and artifacts scattered around. The color palette should be inspired by ancient Egyptian and Mayan cultures, featuring gold, bronze, and jade green tones. The image should be eye-catching and visually appealing, attracting potential players to give the game a try.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Treasures of the Dead and play this exciting slot game for free. Explore ancient treasures and win big!", 2) | Out-Null

# ---------------------------------------------------------------
# 6) Finally, remove the whole "Meta description" paragraph (the
#    paragraph mark too), which was left untouched until now.
# ---------------------------------------------------------------
$d.Paragraphs(2).Range.Delete() | Out-Null

Write-Output "done"
